$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 6091.5713
$ws.Range("I9").Value = 9402.091
$ws.Range("K9").Value = 9402.091
$ws.Range("M9").Value = -9233.091
$ws.Range("H17").Value = 4092176.8
$ws.Range("J17").Value = 4092176.8
$ws.Range("L17").Value = 12276530.4
$ws.Range("N17").Value = -12276866.4
$ws.Range("H18").Value = 439.2857
$ws.Range("I18").Value = 439.2857
$ws.Range("K18").Value = 439.2857
$ws.Range("M18").Value = -155.2857
$ws.Range("H40").Value = 1432710.4
$ws.Range("J40").Value = 1991
$ws.Range("L40").Value = 1991
$ws.Range("N40").Value = -2341
$ws.Range("H44").Value = 3250
$ws.Range("J44").Value = 3250
$ws.Range("L44").Value = 3250
$ws.Range("N44").Value = -4174
$ws.Range("H55").Value = 85.29412000000001
$ws.Range("I55").Value = 74.666664
$ws.Range("J55").Value = 97.25
$ws.Range("K55").Value = 74.666664
$ws.Range("L55").Value = 97.25
$ws.Range("M55").Value = 139.333336
$ws.Range("N55").Value = -525.25
$ws.Range("H70").Value = 3740.4
$ws.Range("I70").Value = 1301
$ws.Range("J70").Value = 5366.6665
$ws.Range("K70").Value = 3903
$ws.Range("L70").Value = 16099.9995
$ws.Range("M70").Value = -3633
$ws.Range("N70").Value = -16639.9995
$ws.Range("H73").Value = 3740.4
$ws.Range("I73").Value = 1301
$ws.Range("J73").Value = 5366.6665
$ws.Range("K73").Value = 3903
$ws.Range("L73").Value = 16099.9995
$ws.Range("M73").Value = -2967
$ws.Range("N73").Value = -17971.9995
$ws.Range("H88").Value = 4110.375
$ws.Range("I88").Value = 3799
$ws.Range("J88").Value = 4629.3335
$ws.Range("K88").Value = 3799
$ws.Range("L88").Value = 4629.3335
$ws.Range("M88").Value = -3393
$ws.Range("N88").Value = -5441.3335
$ws.Range("H91").Value = 4110.375
$ws.Range("I91").Value = 3799
$ws.Range("J91").Value = 4629.3335
$ws.Range("K91").Value = 3799
$ws.Range("L91").Value = 4629.3335
$ws.Range("M91").Value = -2395
$ws.Range("N91").Value = -7437.3335
$ws.Range("H98").Value = 2957.2646
$ws.Range("I98").Value = 2560.9062
$ws.Range("K98").Value = 2560.9062
$ws.Range("M98").Value = -1062.9062
$ws.Range("H107").Value = 1464.7778
$ws.Range("I107").Value = 957.8
$ws.Range("K107").Value = 957.8
$ws.Range("M107").Value = 962.2
$ws.Range("H109").Value = 41488.25
$ws.Range("J109").Value = 38651
$ws.Range("L109").Value = 38651
$ws.Range("N109").Value = -41425
$ws.Range("H112").Value = 2355.75
$ws.Range("I112").Value = 998.8889
$ws.Range("J112").Value = 3465.9092
$ws.Range("K112").Value = 2996.6667
$ws.Range("L112").Value = 10397.7276
$ws.Range("M112").Value = -1888.6667
$ws.Range("N112").Value = -12613.7276
$ws.Range("H122").Value = 2957.2646
$ws.Range("I122").Value = 2560.9062
$ws.Range("K122").Value = 7682.7186
$ws.Range("M122").Value = -5232.7186
$ws.Range("H132").Value = 2947.261
$ws.Range("I132").Value = 3064.2683
$ws.Range("J132").Value = 1987.8
$ws.Range("K132").Value = 9192.804900000001
$ws.Range("L132").Value = 5963.4
$ws.Range("M132").Value = -6662.804900000001
$ws.Range("N132").Value = -11023.4
$ws.Range("H137").Value = 3771.4443
$ws.Range("I137").Value = 1658.5
$ws.Range("K137").Value = 4975.5
$ws.Range("M137").Value = -2425.5
$ws.Range("H138").Value = 2980.1667
$ws.Range("I138").Value = 2319
$ws.Range("J138").Value = 3186.7812
$ws.Range("K138").Value = 6957
$ws.Range("L138").Value = 9560.3436
$ws.Range("M138").Value = -1817
$ws.Range("N138").Value = -19840.3436
$ws.Range("H140").Value = 77712
$ws.Range("J140").Value = 101780
$ws.Range("L140").Value = 101780
$ws.Range("N140").Value = -112140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4049.75
$ws.Range("I2").Value = 3166.3333
$ws.Range("J2").Value = 6700
$ws.Range("K2").Value = 3166.3333
$ws.Range("L2").Value = 6700
$ws.Range("M2").Value = -3053.3333
$ws.Range("N2").Value = -6926
$ws.Range("H32").Value = 1685.622
$ws.Range("I32").Value = 1021.17145
$ws.Range("J32").Value = 5561.5835
$ws.Range("K32").Value = 1021.17145
$ws.Range("L32").Value = 5561.5835
$ws.Range("M32").Value = -734.17145
$ws.Range("N32").Value = -6135.5835
$ws.Range("H45").Value = 1901.2222
$ws.Range("I45").Value = 1901.2222
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1901.2222
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1524.2222
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 45457064
$ws.Range("I61").Value = 45457064
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 45457064
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -45456852
$ws.Range("N61").ClearContents()
$ws.Range("H63").Value = 7345.091
$ws.Range("I63").Value = 2300
$ws.Range("K63").Value = 2300
$ws.Range("M63").Value = -1614
$ws.Range("H66").Value = 7345.091
$ws.Range("I66").Value = 2300
$ws.Range("K66").Value = 11500
$ws.Range("M66").Value = -8068
$ws.Range("H74").Value = 91012456
$ws.Range("I74").Value = 91012456
$ws.Range("K74").Value = 91012456
$ws.Range("M74").Value = -91011582
$ws.Range("H77").Value = 91012456
$ws.Range("I77").Value = 91012456
$ws.Range("K77").Value = 455062280
$ws.Range("M77").Value = -455057912
$ws.Range("H97").Value = 2080.625
$ws.Range("I97").Value = 1879.2858
$ws.Range("J97").Value = 3490
$ws.Range("K97").Value = 1879.2858
$ws.Range("L97").Value = 3490
$ws.Range("M97").Value = -1383.2858
$ws.Range("N97").Value = -4482
$ws.Range("H102").Value = 335422
$ws.Range("I102").Value = 335422
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 335422
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -333800
$ws.Range("N102").ClearContents()
$ws.Range("H110").Value = 15078.434
$ws.Range("I110").Value = 16614.96
$ws.Range("K110").Value = 16614.96
$ws.Range("M110").Value = -14569.96
$ws.Range("H116").Value = 4049.75
$ws.Range("I116").Value = 3166.3333
$ws.Range("J116").Value = 6700
$ws.Range("K116").Value = 3166.3333
$ws.Range("L116").Value = 6700
$ws.Range("M116").Value = -872.3332999999998
$ws.Range("N116").Value = -11288
$ws.Range("H122").Value = 4453.326
$ws.Range("I122").Value = 2179.5
$ws.Range("J122").Value = 5255.853
$ws.Range("K122").Value = 6538.5
$ws.Range("L122").Value = 15767.559
$ws.Range("M122").Value = -4088.5
$ws.Range("N122").Value = -20667.559
$ws.Range("H132").Value = 47700780
$ws.Range("I132").Value = 16650.63
$ws.Range("K132").Value = 49951.89
$ws.Range("M132").Value = -47421.89
$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140
$ws.Range("H136").Value = 45457064
$ws.Range("I136").Value = 45457064
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 136371192
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -136368642
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 49260.332
$ws.Range("J2").Value = 49260.332
$ws.Range("L2").Value = 49260.332
$ws.Range("N2").Value = -49486.332
$ws.Range("H3").Value = 4049.75
$ws.Range("I3").Value = 3166.3333
$ws.Range("J3").Value = 6700
$ws.Range("K3").Value = 3166.3333
$ws.Range("L3").Value = 6700
$ws.Range("M3").Value = -3052.3333
$ws.Range("N3").Value = -6928
$ws.Range("H94").Value = 518.8333
$ws.Range("I94").Value = 541.2963
$ws.Range("J94").Value = 316.66666
$ws.Range("K94").Value = 541.2963
$ws.Range("L94").Value = 316.66666
$ws.Range("M94").Value = -90.29629999999997
$ws.Range("N94").Value = -1218.66666
$ws.Range("H134").Value = 2854.4666
$ws.Range("I134").Value = 2220.3635
$ws.Range("J134").Value = 4598.25
$ws.Range("K134").Value = 6661.0905
$ws.Range("L134").Value = 13794.75
$ws.Range("M134").Value = -4126.0905
$ws.Range("N134").Value = -18864.75
$ws.Range("H141").Value = 80059.664
$ws.Range("J141").Value = 80059.664
$ws.Range("L141").Value = 80059.664
$ws.Range("N141").Value = -90419.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 173344160
$ws.Range("J4").Value = 250016260
$ws.Range("L4").Value = 250016260
$ws.Range("N4").Value = -250016484
$ws.Range("H7").Value = 287.33334
$ws.Range("I7").Value = 191.88889
$ws.Range("J7").Value = 382.77777
$ws.Range("K7").Value = 191.88889
$ws.Range("L7").Value = 382.77777
$ws.Range("M7").Value = -78.88889
$ws.Range("N7").Value = -608.7777699999999
$ws.Range("H16").Value = 2089.75
$ws.Range("I16").Value = 1741.2
$ws.Range("K16").Value = 1741.2
$ws.Range("M16").Value = -1454.2
$ws.Range("H22").Value = 12532.556
$ws.Range("J22").Value = 2388.4
$ws.Range("L22").Value = 2388.4
$ws.Range("N22").Value = -3088.4
$ws.Range("H25").Value = 1570.5
$ws.Range("I25").Value = 134.5
$ws.Range("K25").Value = 134.5
$ws.Range("M25").Value = 39.5
$ws.Range("H31").Value = 2054.027
$ws.Range("I31").Value = 1577.591
$ws.Range("J31").Value = 2752.8
$ws.Range("K31").Value = 1577.591
$ws.Range("L31").Value = 2752.8
$ws.Range("M31").Value = -1282.591
$ws.Range("N31").Value = -3342.8
$ws.Range("H34").Value = 2054.027
$ws.Range("I34").Value = 1577.591
$ws.Range("J34").Value = 2752.8
$ws.Range("K34").Value = 1577.591
$ws.Range("L34").Value = 2752.8
$ws.Range("M34").Value = -1375.591
$ws.Range("N34").Value = -3156.8
$ws.Range("H105").Value = 11353.866
$ws.Range("I105").Value = 5018.364
$ws.Range("K105").Value = 5018.364
$ws.Range("M105").Value = -3271.364
$ws.Range("H113").Value = 2089.75
$ws.Range("I113").Value = 1741.2
$ws.Range("K113").Value = 1741.2
$ws.Range("M113").Value = 428.8
$ws.Range("H122").Value = 3653
$ws.Range("I122").Value = 3670.6667
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 11012.0001
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -8562.000100000001
$ws.Range("N122").Value = -15700
$ws.Range("H132").Value = 41754.8
$ws.Range("I132").Value = 50389.047
$ws.Range("J132").Value = 2421
$ws.Range("K132").Value = 151167.141
$ws.Range("L132").Value = 7263
$ws.Range("M132").Value = -148637.141
$ws.Range("N132").Value = -12323
$ws.Range("H134").Value = 2645.6316
$ws.Range("I134").Value = 2667.8667
$ws.Range("J134").Value = 2562.25
$ws.Range("K134").Value = 8003.6001
$ws.Range("L134").Value = 7686.75
$ws.Range("M134").Value = -5468.6001
$ws.Range("N134").Value = -12756.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4238.9
$ws.Range("J3").Value = 10500
$ws.Range("L3").Value = 31500
$ws.Range("N3").Value = -31724
$ws.Range("H12").Value = 45.666668
$ws.Range("J12").Value = 47.25
$ws.Range("L12").Value = 141.75
$ws.Range("N12").Value = -487.75
$ws.Range("H56").Value = 15382.8
$ws.Range("I56").Value = 15382.8
$ws.Range("K56").Value = 15382.8
$ws.Range("M56").Value = -14852.8
$ws.Range("H132").Value = 4353.8887
$ws.Range("I132").Value = 2563.182
$ws.Range("J132").Value = 7167.857
$ws.Range("K132").Value = 23068.638
$ws.Range("L132").Value = 64510.713
$ws.Range("M132").Value = -20538.638
$ws.Range("N132").Value = -69570.713
$ws.Range("H133").Value = 2000
$ws.Range("I133").Value = 1000
$ws.Range("J133").Value = 3000
$ws.Range("K133").Value = 3000
$ws.Range("L133").Value = 9000
$ws.Range("N133").Value = -19120
$ws.Range("M133").Value = 2060
$ws.Range("H134").Value = 2590.3438
$ws.Range("I134").Value = 559.72
$ws.Range("K134").Value = 1679.16
$ws.Range("M134").Value = 3390.84
$ws.Range("H138").Value = 3611.7778
$ws.Range("I138").Value = 3684.125
$ws.Range("K138").Value = 11052.375
$ws.Range("M138").Value = -5912.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1563265.8
$ws.Range("I2").Value = 2500078.2
$ws.Range("K2").Value = 2500078.2
$ws.Range("M2").Value = -2499965.2
$ws.Range("H5").Value = 28999.5
$ws.Range("I5").Value = 8000
$ws.Range("K5").Value = 8000
$ws.Range("M5").Value = -7888
$ws.Range("H15").Value = 61746.75
$ws.Range("J15").Value = 61746.75
$ws.Range("L15").Value = 61746.75
$ws.Range("N15").Value = -62322.75
$ws.Range("H80").Value = 5718.1816
$ws.Range("I80").Value = 6276.091
$ws.Range("K80").Value = 6276.091
$ws.Range("M80").Value = -5278.091
$ws.Range("H81").Value = 61746.75
$ws.Range("J81").Value = 61746.75
$ws.Range("L81").Value = 61746.75
$ws.Range("N81").Value = -63742.75
$ws.Range("H83").Value = 5718.1816
$ws.Range("I83").Value = 6276.091
$ws.Range("K83").Value = 31380.455
$ws.Range("M83").Value = -26388.455
$ws.Range("H84").Value = 61746.75
$ws.Range("J84").Value = 61746.75
$ws.Range("L84").Value = 185240.25
$ws.Range("N84").Value = -195224.25
$ws.Range("H113").Value = 4856.857
$ws.Range("I113").Value = 4499
$ws.Range("K113").Value = 4499
$ws.Range("M113").Value = -2329
$ws.Range("H122").Value = 3141.739
$ws.Range("I122").Value = 2995.75
$ws.Range("J122").Value = 3475.4285
$ws.Range("K122").Value = 8987.25
$ws.Range("L122").Value = 10426.2855
$ws.Range("M122").Value = -6537.25
$ws.Range("N122").Value = -15326.2855
$ws.Range("H132").Value = 7399.7905
$ws.Range("I132").Value = 5224.4585
$ws.Range("J132").Value = 14858.071
$ws.Range("K132").Value = 15673.3755
$ws.Range("L132").Value = 44574.213
$ws.Range("M132").Value = -13143.3755
$ws.Range("N132").Value = -49634.213

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6074.857
$ws.Range("I7").Value = 4021.4285
$ws.Range("J7").Value = 8128.2856
$ws.Range("K7").Value = 4021.4285
$ws.Range("L7").Value = 8128.2856
$ws.Range("M7").Value = -3909.4285
$ws.Range("N7").Value = -8352.285599999999
$ws.Range("H22").Value = 2264.9
$ws.Range("I22").Value = 1665.5714
$ws.Range("J22").Value = 3663.3333
$ws.Range("K22").Value = 1665.5714
$ws.Range("L22").Value = 3663.3333
$ws.Range("M22").Value = -1370.5714
$ws.Range("N22").Value = -4253.3333
$ws.Range("H27").Value = 2264.9
$ws.Range("I27").Value = 1665.5714
$ws.Range("J27").Value = 3663.3333
$ws.Range("K27").Value = 1665.5714
$ws.Range("L27").Value = 3663.3333
$ws.Range("M27").Value = -1558.5714
$ws.Range("N27").Value = -3877.3333
$ws.Range("H40").Value = 8419.615
$ws.Range("I40").Value = 8445.1
$ws.Range("K40").Value = 8445.1
$ws.Range("M40").Value = -8309.1
$ws.Range("H46").Value = 1740.7142
$ws.Range("I46").Value = 988.1818
$ws.Range("J46").Value = 4500
$ws.Range("K46").Value = 988.1818
$ws.Range("L46").Value = 4500
$ws.Range("M46").Value = -800.1818
$ws.Range("N46").Value = -4876
$ws.Range("H55").Value = 650
$ws.Range("I55").Value = 300.5
$ws.Range("J55").Value = 999.5
$ws.Range("K55").Value = 300.5
$ws.Range("L55").Value = 999.5
$ws.Range("N55").Value = -1345.5
$ws.Range("M55").Value = -127.5
$ws.Range("H61").Value = 2720.875
$ws.Range("I61").Value = 2720.875
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2720.875
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2518.875
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 8749.25
$ws.Range("I68").Value = 8998.5
$ws.Range("J68").Value = 8500
$ws.Range("K68").Value = 8998.5
$ws.Range("L68").Value = 8500
$ws.Range("N68").Value = -9998
$ws.Range("M68").Value = -8249.5
$ws.Range("H71").Value = 8749.25
$ws.Range("I71").Value = 8998.5
$ws.Range("J71").Value = 8500
$ws.Range("K71").Value = 44992.5
$ws.Range("L71").Value = 42500
$ws.Range("N71").Value = -49988
$ws.Range("M71").Value = -41248.5
$ws.Range("H113").Value = 2720.875
$ws.Range("I113").Value = 2720.875
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2720.875
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -550.875
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 4225.85
$ws.Range("I122").Value = 3643.6875
$ws.Range("J122").Value = 6554.5
$ws.Range("K122").Value = 10931.0625
$ws.Range("L122").Value = 19663.5
$ws.Range("M122").Value = -8481.0625
$ws.Range("N122").Value = -24563.5
$ws.Range("H126").Value = 6074.857
$ws.Range("I126").Value = 4021.4285
$ws.Range("J126").Value = 8128.2856
$ws.Range("K126").Value = 12064.2855
$ws.Range("L126").Value = 24384.8568
$ws.Range("M126").Value = -9594.2855
$ws.Range("N126").Value = -29324.8568
$ws.Range("H132").Value = 3098
$ws.Range("I132").Value = 2990.7856
$ws.Range("K132").Value = 8972.356800000001
$ws.Range("M132").Value = -6442.356800000001
$ws.Range("H136").Value = 692425.25
$ws.Range("I136").Value = 1054728.4
$ws.Range("J136").Value = 4049.4
$ws.Range("K136").Value = 3164185.2
$ws.Range("L136").Value = 12148.2
$ws.Range("M136").Value = -3161635.2
$ws.Range("N136").Value = -17248.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3854.4333
$ws.Range("I81").Value = 2557.7
$ws.Range("J81").Value = 6447.9
$ws.Range("K81").Value = 5115.4
$ws.Range("L81").Value = 12895.8
$ws.Range("M81").Value = -4054.4
$ws.Range("N81").Value = -15017.8
$ws.Range("H84").Value = 3854.4333
$ws.Range("I84").Value = 2557.7
$ws.Range("J84").Value = 6447.9
$ws.Range("K84").Value = 25577
$ws.Range("L84").Value = 64479
$ws.Range("M84").Value = -20273
$ws.Range("N84").Value = -75087
$ws.Range("H100").Value = 1251749.8
$ws.Range("I100").Value = 2001199.8
$ws.Range("K100").Value = 4002399.6
$ws.Range("M100").Value = -4001858.6
$ws.Range("H104").Value = 24267
$ws.Range("J104").Value = 24267
$ws.Range("L104").Value = 24267
$ws.Range("N104").Value = -31255
$ws.Range("H107").Value = 1302.0303
$ws.Range("I107").Value = 1006.85
$ws.Range("J107").Value = 1756.1538
$ws.Range("K107").Value = 3020.55
$ws.Range("L107").Value = 5268.4614
$ws.Range("M107").Value = -1100.55
$ws.Range("N107").Value = -9108.4614
$ws.Range("H113").Value = 834.6111
$ws.Range("I113").Value = 795.9167
$ws.Range("J113").Value = 912
$ws.Range("K113").Value = 2387.7501
$ws.Range("L113").Value = 2736
$ws.Range("M113").Value = -217.7501000000002
$ws.Range("N113").Value = -7076
$ws.Range("H122").Value = 87044.75
$ws.Range("I122").Value = 103440.5
$ws.Range("J122").Value = 5066
$ws.Range("K122").Value = 310321.5
$ws.Range("L122").Value = 15198
$ws.Range("M122").Value = -307871.5
$ws.Range("N122").Value = -20098
$ws.Range("H126").Value = 7399.8
$ws.Range("I126").Value = 6249.75
$ws.Range("J126").Value = 12000
$ws.Range("K126").Value = 18749.25
$ws.Range("L126").Value = 36000
$ws.Range("M126").Value = -16279.25
$ws.Range("N126").Value = -40940
$ws.Range("H132").Value = 2072.3794
$ws.Range("I132").Value = 1931.1305
$ws.Range("J132").Value = 2613.8333
$ws.Range("K132").Value = 5793.3915
$ws.Range("L132").Value = 7841.499899999999
$ws.Range("M132").Value = -3263.3915
$ws.Range("N132").Value = -12901.4999
$ws.Range("H136").Value = 1821.5161
$ws.Range("I136").Value = 1659
$ws.Range("J136").Value = 2218.7778
$ws.Range("K136").Value = 4977
$ws.Range("L136").Value = 6656.3334
$ws.Range("M136").Value = -2427
$ws.Range("N136").Value = -11756.3334

